# loading_percent.xlsx - Case_1_46 (380 kV case) results refresh.
# Re-writes the simulated loading-percent results in columns B:M (rows 2-25,
# one row per line case) with the newly computed values. Columns A, N and O
# (case index / always-zero columns) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (case 0): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.038478590194587
$row[0,1] = 0
$row[0,2] = 4.200236232346084
$row[0,3] = 10.58592210995402
$row[0,4] = 55.97689310170403
$row[0,5] = 3.775357352031684
$row[0,6] = 0
$row[0,7] = 28.56943258219281
$row[0,8] = 10.21095220581566
$row[0,9] = 19.77108203066788
$row[0,10] = 0
$row[0,11] = 19.98108532810528
$ws.Range("B2:M2").Value = $row

# Row 3 (case 1): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 4.971485998590704
$row[0,1] = 0
$row[0,2] = 4.206681548598971
$row[0,3] = 10.60911583372151
$row[0,4] = 55.66012800396597
$row[0,5] = 3.779547390959261
$row[0,6] = 0
$row[0,7] = 28.48289790478313
$row[0,8] = 10.22610300432585
$row[0,9] = 19.62194547617872
$row[0,10] = 0
$row[0,11] = 19.95747371372629
$ws.Range("B3:M3").Value = $row

# Row 4 (case 2): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 4.928570664651845
$row[0,1] = 0
$row[0,2] = 4.21112857077487
$row[0,3] = 10.62478750365319
$row[0,4] = 55.47508247768989
$row[0,5] = 3.782251838293728
$row[0,6] = 0
$row[0,7] = 28.43198611629305
$row[0,8] = 10.23640394129876
$row[0,9] = 19.53621358696942
$row[0,10] = 0
$row[0,11] = 19.94799752428383
$ws.Range("B4:M4").Value = $row

# Row 5 (case 3): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 4.910639977739941
$row[0,1] = 0
$row[0,2] = 4.213063797015367
$row[0,3] = 10.63153387264049
$row[0,4] = 55.40209289360868
$row[0,5] = 3.78338718399052
$row[0,6] = 0
$row[0,7] = 28.41180174333098
$row[0,8] = 10.24085281191987
$row[0,9] = 19.50277854823541
$row[0,10] = 0
$row[0,11] = 19.94540092869707
$ws.Range("B5:M5").Value = $row

# Row 6 (case 4): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 4.907636055765443
$row[0,1] = 0
$row[0,2] = 4.21339256814675
$row[0,3] = 10.63267585581099
$row[0,4] = 55.39012021583514
$row[0,5] = 3.783577720082534
$row[0,6] = 0
$row[0,7] = 28.40848423343001
$row[0,8] = 10.24160671631088
$row[0,9] = 19.49731827885153
$row[0,10] = 0
$row[0,11] = 19.94504621969121
$ws.Range("B6:M6").Value = $row

# Row 7 (case 5): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 4.928330628494947
$row[0,1] = 0
$row[0,2] = 4.211154171864192
$row[0,3] = 10.62487702938562
$row[0,4] = 55.47408827224777
$row[0,5] = 3.782267015115165
$row[0,6] = 0
$row[0,7] = 28.43171161989736
$row[0,8] = 10.236462923244
$row[0,9] = 19.53575654963775
$row[0,10] = 0
$row[0,11] = 19.94795738135101
$ws.Range("B7:M7").Value = $row

# Row 8 (case 6): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.015752618043939
$row[0,1] = 0
$row[0,2] = 4.202356920233945
$row[0,3] = 10.59362255663338
$row[0,4] = 55.86573743273978
$row[0,5] = 3.776774813550891
$row[0,6] = 0
$row[0,7] = 28.53913711265431
$row[0,8] = 10.21596908506249
$row[0,9] = 19.71846880783853
$row[0,10] = 0
$row[0,11] = 19.97190354267261
$ws.Range("B8:M8").Value = $row

# Row 9 (case 7): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.172818675410919
$row[0,1] = 0
$row[0,2] = 4.188993756282844
$row[0,3] = 10.54367121793453
$row[0,4] = 56.70674065828276
$row[0,5] = 3.767043915417692
$row[0,6] = 0
$row[0,7] = 28.76722206447828
$row[0,8] = 10.1836962587671
$row[0,9] = 20.12149261533284
$row[0,10] = 0
$row[0,11] = 20.05856008849599
$ws.Range("B9:M9").Value = $row

# Row 10 (case 8): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.279168700597397
$row[0,1] = 0
$row[0,2] = 4.181551524709167
$row[0,3] = 10.51386708899135
$row[0,4] = 57.36622632168421
$row[0,5] = 3.760519723704313
$row[0,6] = 0
$row[0,7] = 28.94513158071103
$row[0,8] = 10.16480395835606
$row[0,9] = 20.44253718826056
$row[0,10] = 0
$row[0,11] = 20.14617354301249
$ws.Range("B10:M10").Value = $row

# Row 11 (case 9): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.325530854144735
$row[0,1] = 0
$row[0,2] = 4.178682825672948
$row[0,3] = 10.50180185096604
$row[0,4] = 57.67459768317192
$row[0,5] = 3.757685619164441
$row[0,6] = 0
$row[0,7] = 29.02824044366515
$row[0,8] = 10.1572541252913
$row[0,9] = 20.59346571113637
$row[0,10] = 0
$row[0,11] = 20.19115608210011
$ws.Range("B11:M11").Value = $row

# Row 12 (case 10): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.342792654225226
$row[0,1] = 0
$row[0,2] = 4.177670920938255
$row[0,3] = 10.49744740049709
$row[0,4] = 57.79251160030181
$row[0,5] = 3.756631515338199
$row[0,6] = 0
$row[0,7] = 29.06001758396022
$row[0,8] = 10.15454524089844
$row[0,9] = 20.65127302334923
$row[0,10] = 0
$row[0,11] = 20.20891901793821
$ws.Range("B12:M12").Value = $row

# Row 13 (case 11): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.339088185124014
$row[0,1] = 0
$row[0,2] = 4.177885541927023
$row[0,3] = 10.49837567855359
$row[0,4] = 57.76706703563971
$row[0,5] = 3.756857687563049
$row[0,6] = 0
$row[0,7] = 29.05316033896061
$row[0,8] = 10.15512197481361
$row[0,9] = 20.63879488096435
$row[0,10] = 0
$row[0,11] = 20.2050611690971
$ws.Range("B13:M13").Value = $row

# Row 14 (case 12): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.326956920630791
$row[0,1] = 0
$row[0,2] = 4.178598083870502
$row[0,3] = 10.50143931253585
$row[0,4] = 57.6842760457171
$row[0,5] = 3.757598515145581
$row[0,6] = 0
$row[0,7] = 29.03084868728402
$row[0,8] = 10.15702825667419
$row[0,9] = 20.5982087281184
$row[0,10] = 0
$row[0,11] = 20.19260288917422
$ws.Range("B14:M14").Value = $row

# Row 15 (case 13): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.319487684218261
$row[0,1] = 0
$row[0,2] = 4.179044229358999
$row[0,3] = 10.50334378891577
$row[0,4] = 57.63371077812678
$row[0,5] = 3.758054778336334
$row[0,6] = 0
$row[0,7] = 29.01722172847516
$row[0,8] = 10.1582154501105
$row[0,9] = 20.57343224797989
$row[0,10] = 0
$row[0,11] = 20.18506651371328
$ws.Range("B15:M15").Value = $row

# Row 16 (case 14): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.276097846482042
$row[0,1] = 0
$row[0,2] = 4.181749408493574
$row[0,3] = 10.514685596922
$row[0,4] = 57.34623658638864
$row[0,5] = 3.760707620227001
$row[0,6] = 0
$row[0,7] = 28.93974340283465
$row[0,8] = 10.1653183642227
$row[0,9] = 20.43276782500462
$row[0,10] = 0
$row[0,11] = 20.14333631322263
$ws.Range("B16:M16").Value = $row

# Row 17 (case 15): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.248959691677919
$row[0,1] = 0
$row[0,2] = 4.18354138292962
$row[0,3] = 10.5220255839791
$row[0,4] = 57.17197969400149
$row[0,5] = 3.762369226806924
$row[0,6] = 0
$row[0,7] = 28.89276539444033
$row[0,8] = 10.16994318710516
$row[0,9] = 20.34769162124374
$row[0,10] = 0
$row[0,11] = 20.11904374867688
$ws.Range("B17:M17").Value = $row

# Row 18 (case 16): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.23316094326189
$row[0,1] = 0
$row[0,2] = 4.184620718512607
$row[0,3] = 10.52638787765896
$row[0,4] = 57.07254330020112
$row[0,5] = 3.763337538872963
$row[0,6] = 0
$row[0,7] = 28.86595026474626
$row[0,8] = 10.172701570646
$row[0,9] = 20.29922070130825
$row[0,10] = 0
$row[0,11] = 20.10555441457969
$ws.Range("B18:M18").Value = $row

# Row 19 (case 17): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.227779336334278
$row[0,1] = 0
$row[0,2] = 4.184994513944422
$row[0,3] = 10.52788901909415
$row[0,4] = 57.03901367688671
$row[0,5] = 3.763667560644307
$row[0,6] = 0
$row[0,7] = 28.8569066147896
$row[0,8] = 10.17365239945617
$row[0,9] = 20.28289019519197
$row[0,10] = 0
$row[0,11] = 20.10107036542049
$ws.Range("B19:M19").Value = $row

# Row 20 (case 18): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.251868250901742
$row[0,1] = 0
$row[0,2] = 4.183345589525899
$row[0,3] = 10.5212296877186
$row[0,4] = 57.19044816272255
$row[0,5] = 3.762191042834331
$row[0,6] = 0
$row[0,7] = 28.89774508003498
$row[0,8] = 10.16944069308094
$row[0,9] = 20.35670059078091
$row[0,10] = 0
$row[0,11] = 20.12157979022549
$ws.Range("B20:M20").Value = $row

# Row 21 (case 19): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.330528195283616
$row[0,1] = 0
$row[0,2] = 4.178386772864575
$row[0,3] = 10.50053363261483
$row[0,4] = 57.70856332949759
$row[0,5] = 3.757380398501946
$row[0,6] = 0
$row[0,7] = 29.03739392732984
$row[0,8] = 10.15646426340545
$row[0,9] = 20.61011250764423
$row[0,10] = 0
$row[0,11] = 20.19624247257734
$ws.Range("B21:M21").Value = $row

# Row 22 (case 20): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.380218316366925
$row[0,1] = 0
$row[0,2] = 4.175579654057779
$row[0,3] = 10.48825707871434
$row[0,4] = 58.05379782568928
$row[0,5] = 3.754347693295238
$row[0,6] = 0
$row[0,7] = 29.13044012481618
$row[0,8] = 10.14885810596833
$row[0,9] = 20.77952045003238
$row[0,10] = 0
$row[0,11] = 20.24928368526291
$ws.Range("B22:M22").Value = $row

# Row 23 (case 21): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.353856437773099
$row[0,1] = 0
$row[0,2] = 4.177038148088041
$row[0,3] = 10.49469506836597
$row[0,4] = 57.86895563907122
$row[0,5] = 3.755956160661948
$row[0,6] = 0
$row[0,7] = 29.08061953981694
$row[0,8] = 10.15283765784025
$row[0,9] = 20.68877382937696
$row[0,10] = 0
$row[0,11] = 20.22058914548418
$ws.Range("B23:M23").Value = $row

# Row 24 (case 22): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.250553903489874
$row[0,1] = 0
$row[0,2] = 4.183433954777735
$row[0,3] = 10.52158906873695
$row[0,4] = 57.18209623897542
$row[0,5] = 3.762271559205928
$row[0,6] = 0
$row[0,7] = 28.89549316199145
$row[0,8] = 10.16966756066601
$row[0,9] = 20.35262625987218
$row[0,10] = 0
$row[0,11] = 20.12043176015311
$ws.Range("B24:M24").Value = $row

# Row 25 (case 23): columns B:M
$row = New-Object 'object[,]' 1,12
$row[0,0] = 5.131897147806348
$row[0,1] = 0
$row[0,2] = 4.192192020818051
$row[0,3] = 10.55597224274384
$row[0,4] = 56.47169954061116
$row[0,5] = 3.769565995937047
$row[0,6] = 0
$row[0,7] = 28.70368120063769
$row[0,8] = 10.19158019798304
$row[0,9] = 20.00790656397835
$row[0,10] = 0
$row[0,11] = 20.03088944481266
$ws.Range("B25:M25").Value = $row
